$d = $word.ActiveDocument

# 1) Merge "Problem 1. Problem 2" into a single bold run "Problem 1, Problem 2"
$range = $d.Content
$range.Find.Execute("Problem 1. Problem 2", $false, $false, $false, $false, $false, $true, 1, $false, "Problem 1, Problem 2", 2)

# 2) Update the trailing sentence wording
$range2 = $d.Content
$range2.Find.Execute("It is larger than others and can" + [char]8217 + "t be took as an optimal plan. In therms of Plan lenght we can say that all other algorithms tested offer optimal plans.", $false, $false, $false, $false, $false, $true, 1, $false, "It is larger than others and doesn" + [char]8217 + "t seem to be an optimal plan. In therms of Plan lenght we can say that other algorithms tested offer optimal plans.", 2)
